$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend formatting for the two new rows (23, 24) by copying the
# format from the existing last data row (22), so the new cells pick up
# style index 2 just like the rest of the table.
$ws.Range("A22:C22").Copy()
$ws.Range("A23:C24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 21: fill in Description (B21), previously empty
$ws.Range("B21").Value = "Done"

# Row 22: fill in Description (B22), previously empty, and change
# Runmode (C22) from Y to N
$ws.Range("B22").Value = "Changes where the event replay or not"
$ws.Range("C22").Value = "N"

# Row 23 (new): ReplayReport
$ws.Range("A23").Value = "ReplayReport"
$ws.Range("B23").Value = "Changes where the report replay or not"
$ws.Range("C23").Value = "N"

# Row 24 (new): ImportSIMs
$ws.Range("A24").Value = "ImportSIMs"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "Y"

# Update the active selection to match the diff (B20)
$ws.Range("B20").Select()
